$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.023.43"
$ws.Range("E2").Value = "  -7.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.903.66"
$ws.Range("E3").Value = "  -9.67%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.99"
$ws.Range("E5").Value = "  -11.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.70"
$ws.Range("E6").Value = "  -16.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.874.68"
$ws.Range("E8").Value = "  -10.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.451"
$ws.Range("E9").Value = "  -17.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  -19.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.66"
$ws.Range("E11").Value = "  -13.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  -15.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.97"
$ws.Range("E13").Value = "  -20.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("E14").Value = "  -20.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.395.53"
$ws.Range("E15").Value = "  -9.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.952.15"
$ws.Range("E16").Value = "  -7.36%  "
$ws.Range("E17").Value = "  -5.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.914.65"
$ws.Range("E18").Value = "  -9.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "461.20"
$ws.Range("E19").Value = "  -13.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("E20").Value = "  -15.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.30"
$ws.Range("E21").Value = "  -17.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.616"
$ws.Range("E22").Value = "  -19.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.25"
$ws.Range("E23").Value = "  -21.49%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.72"
$ws.Range("E24").Value = "  -15.76%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.52"
$ws.Range("E25").Value = "  -14.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -19.16%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("E28").Value = "  -16.48%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.78"
$ws.Range("E29").Value = "  -18.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.92"
$ws.Range("E30").Value = "  -18.19%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.03"
$ws.Range("E32").Value = "  -11.30%  "
$ws.Range("E33").Value = "  -15.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.85"
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "446.21"
$ws.Range("E35").Value = "  -18.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.26"
$ws.Range("E36").Value = "  -19.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.49"
$ws.Range("E37").Value = "  -21.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0372"
$ws.Range("E38").Value = "  -13.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0724"
$ws.Range("E39").Value = "  -16.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("E40").Value = "  -11.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.52"
$ws.Range("E41").Value = "  -19.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.587.52"
$ws.Range("E42").Value = "  -11.11%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("E44").Value = "  -22.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.211"
$ws.Range("E45").Value = "  -19.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.07"
$ws.Range("E46").Value = "  -8.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0981"
$ws.Range("E47").Value = "  -14.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.74"
$ws.Range("E48").Value = "  -19.72%  "
$ws.Range("B49").Value = "BitgetToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = [string]::Concat("0.0", $sub3, "0444")
$ws.Range("E50").Value = "  -24.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.67"
$ws.Range("E51").Value = "  -22.19%  "
